$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.600.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.524.36'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.522.30'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +3.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('E11').Value = '  -4.88%  '
$ws.Range('E12').Value = '  -2.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.117.18'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.98%  '
$ws.Range('E14').Value = '  -6.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.67'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.534.68'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.412.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.79%  '
$ws.Range('E21').Value = '  -3.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '422.80'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.589'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.670.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.66%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000113'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.86%  '
$ws.Range('E28').Value = '  -2.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.23%  '
$ws.Range('E30').Value = '  -6.73%  '
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.531.74'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.58%  '
$ws.Range('E33').Value = '  -2.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.94%  '
$ws.Range('E36').Value = '  -9.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.54'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.62'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '173.29'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('E40').Value = '  -8.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0807'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.97'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.54%  '
$ws.Range('E43').Value = '  -5.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.32%  '
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.900'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.25%  '
